$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (values that Excel will not reinterpret as numbers)
$textUpdates = @{
    "D2" = "26.116.18"
    "E2" = "  -0.28%  "
    "D3" = "1.653.38"
    "E3" = "  -0.37%  "
    "E4" = "  -0.33%  "
    "E5" = "  +0.17%  "
    "E6" = "  +1.60%  "
    "E7" = "  -0.28%  "
    "E8" = "  -1.99%  "
    "E9" = "  +0.24%  "
    "E10" = "  -2.82%  "
    "E11" = "  +0.26%  "
    "E12" = "  +1.37%  "
    "D13" = "1.640.21"
    "E13" = "  -1.27%  "
    "E14" = "  +0.23%  "
    "D15" = "0.0₅8132"
    "E15" = "  -0.58%  "
    "E16" = "  +1.02%  "
    "D17" = "26.138.57"
    "E17" = "  -0.30%  "
    "E18" = "  -0.29%  "
    "E19" = "  -2.41%  "
    "E20" = "  +1.00%  "
    "E21" = "  -0.88%  "
    "E22" = "  -1.29%  "
    "E23" = "  -0.42%  "
    "E24" = "  +1.11%  "
    "E25" = "  +0.58%  "
    "E26" = "  +0.86%  "
    "E27" = "  -0.08%  "
    "E28" = "  +1.74%  "
    "E29" = "  -1.02%  "
    "E30" = "  -0.24%  "
    "E31" = "  -2.14%  "
    "E32" = "  -2.36%  "
    "E33" = "  -5.52%  "
    "E34" = "  +0.04%  "
    "E35" = "  -3.22%  "
    "E36" = "  -0.81%  "
    "E38" = "  +1.16%  "
    "E39" = "  -1.44%  "
    "E40" = "  -1.92%  "
    "E41" = "  -0.18%  "
    "B42" = "Quant"
    "C42" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
    "E42" = "  +1.49%  "
    "B43" = "Maker"
    "C43" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D43" = "1.007.72"
    "E43" = "  -2.30%  "
    "D44" = "1.797.93"
    "E44" = "  -0.15%  "
    "E45" = "  -0.16%  "
    "E46" = "  +0.48%  "
    "E47" = "  -0.38%  "
    "E48" = "  +1.39%  "
    "E49" = "  +1.14%  "
    "E50" = "  -0.63%  "
    "E51" = "  -4.06%  "
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Numeric-looking price values that must stay as literal text strings
# (force text format so Excel does not coerce them into numbers, then
#  strip the format change so no stray style is left behind)
$forcedTextUpdates = @{
    "D5" = "218.38"
    "D6" = "0.5289"
    "D7" = "1.002"
    "D10" = "20.38"
    "D11" = "0.07742"
    "D12" = "4.478"
    "D14" = "0.5461"
    "D16" = "65.28"
    "D20" = "194.11"
    "D21" = "10.03"
    "D22" = "5.984"
    "D24" = "140.06"
    "D26" = "7.261"
    "D27" = "16.16"
    "D28" = "1.437"
    "D29" = "0.05925"
    "D30" = "1.278"
    "D32" = "3.230"
    "D33" = "1.544"
    "D34" = "2.411"
    "D35" = "0.9446"
    "D36" = "2.758"
    "D37" = "0.5635"
    "D39" = "5.851"
    "D40" = "0.8438"
    "D42" = "100.93"
    "D45" = "56.82"
    "D47" = "1.005"
    "D48" = "0.4288"
    "D49" = "1.472"
    "D50" = "0.05151"
    "D51" = "7.745"
}
foreach ($ref in $forcedTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $forcedTextUpdates[$ref]
    $cell.ClearFormats()
}

